$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("I2").Value = 0.8652647215324584
$ws.Range("J2").Value = 0.8652647215324583
$ws.Range("M2").Value = 547.450775
$ws.Range("N2").Value = 1642.352325
$ws.Range("O2").Value = 0.8253533007282613
$ws.Range("P2").Value = 0.8253533007282614
$ws.Range("Q2").Value = 71.21768883615833
$ws.Range("R2").Value = 640.959199525425
$ws.Range("S2").Value = 0.7141490939205344
$ws.Range("T2").Value = 0.7141490939205344

# Row 3
$ws.Range("I3").Value = 0.8652647215324584
$ws.Range("J3").Value = 0.8652647215324583
$ws.Range("O3").Value = 0.002183077622430991
$ws.Range("P3").Value = 0.002183077622430991
$ws.Range("S3").Value = 0.001888940051056493
$ws.Range("T3").Value = 0.001888940051056493

# Row 4
$ws.Range("I4").Value = 0.8652647215324584
$ws.Range("J4").Value = 0.8652647215324583
$ws.Range("M4").Value = 114.393852
$ws.Range("N4").Value = 343.181556
$ws.Range("O4").Value = 0.1724636216493076
$ws.Range("P4").Value = 0.1724636216493076
$ws.Range("Q4").Value = 14.881458075396
$ws.Range("R4").Value = 133.933122678564
$ws.Range("S4").Value = 0.1492266875608674
$ws.Range("T4").Value = 0.1492266875608674

# Row 5
$ws.Range("E5").Value = 2
$ws.Range("F5").Value = 0.6666666666666666
$ws.Range("G5").Value = 0.020257
$ws.Range("H5").Value = 0.060771
$ws.Range("I5").Value = 0.1347352784675417
$ws.Range("J5").Value = 0.1347352784675417
$ws.Range("M5").Value = 547.450775
$ws.Range("N5").Value = 1642.352325
$ws.Range("O5").Value = 0.8253533007282613
$ws.Range("P5").Value = 0.8253533007282614
$ws.Range("Q5").Value = 11.089710349175
$ws.Range("R5").Value = 99.80739314257499
$ws.Range("S5").Value = 0.111204206807727
$ws.Range("T5").Value = 0.111204206807727

# Row 6
$ws.Range("E6").Value = 2
$ws.Range("F6").Value = 0.6666666666666666
$ws.Range("G6").Value = 0.020257
$ws.Range("H6").Value = 0.060771
$ws.Range("I6").Value = 0.1347352784675417
$ws.Range("J6").Value = 0.1347352784675417
$ws.Range("O6").Value = 0.002183077622430991
$ws.Range("P6").Value = 0.002183077622430991
$ws.Range("Q6").Value = 0.02933252763533334
$ws.Range("R6").Value = 0.263992748718
$ws.Range("S6").Value = 0.0002941375713744984
$ws.Range("T6").Value = 0.0002941375713744984

# Row 7
$ws.Range("E7").Value = 2
$ws.Range("F7").Value = 0.6666666666666666
$ws.Range("G7").Value = 0.020257
$ws.Range("H7").Value = 0.060771
$ws.Range("I7").Value = 0.1347352784675417
$ws.Range("J7").Value = 0.1347352784675417
$ws.Range("M7").Value = 114.393852
$ws.Range("N7").Value = 343.181556
$ws.Range("O7").Value = 0.1724636216493076
$ws.Range("P7").Value = 0.1724636216493076
$ws.Range("Q7").Value = 2.317276259964
$ws.Range("R7").Value = 20.855486339676
$ws.Range("S7").Value = 0.0232369340884402
$ws.Range("T7").Value = 0.02323693408844021
